$d = $word.ActiveDocument

# 1. Remove the "Meta description: Read our review..." paragraph that
#    currently sits right after the H1 title paragraph.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. At the end of the document, the paragraph that used to hold the
#    italic image-prompt text is now split into two paragraphs:
#      - a new bold paragraph with the page title text
#      - the original paragraph (still italic) but with new text
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

# Insert a new paragraph break + placeholder text right before the
# existing (italic) text, without inheriting the italic formatting
# (Range.Text assignment at a zero-length range does not copy rPr).
$insertPos = $lastPara.Range.Start
$insRange = $d.Range($insertPos, $insertPos)
$insRange.Text = "Play Cleopatra Diamond Spins Free - A Unique Ancient Egypt Theme`r"

# Make the freshly inserted paragraph's text bold (exclude the paragraph
# mark itself so no extra pPr/rPr gets attached to it).
$newParaIndex = $count
$newPara = $d.Paragraphs($newParaIndex)
$newTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newTextRange.Font.Bold = 1

# 3. Replace the text of the (now shifted) final paragraph - keep its
#    existing italic run formatting, just swap the wording.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalTextRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalTextRange.Text = "Read our review on Cleopatra Diamond Spins - an online slot game with unique Ancient Egypt themes. Play for free today."
